# practice: Solve Leetcode Problem# 176 Second Highest Salary
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (Anna's entry for "176. Second Highest Salary") previously only had
# A:D filled in. Add the start-date placeholder in column E, matching the
# pattern used by the other "Anna" rows (e.g. row 21).
$ws.Range("E23").Value = "2020/12/"

# Row 24 (Stephan's entry for the same problem) gets the completion date,
# category and status filled in now that the problem has been solved.
$ws.Range("E24").Value = "2020/12/14"
$ws.Range("F24").Value = "SQL"
$ws.Range("G24").Value = "Completed"

# A new problem row pair is started: "Anna" row 25 and "Stephan" row 26,
# both currently Easy / LeetCode / not yet started (only A:C filled in).
# Columns A & B carry the same left-aligned style used throughout the rest
# of the table (e.g. A2:B24), column C keeps the default (unstyled) look.
$ws.Range("A25").Value = "LeetCode"
$ws.Range("A25").HorizontalAlignment = -4131
$ws.Range("B25").Value = "Anna"
$ws.Range("B25").HorizontalAlignment = -4131
$ws.Range("C25").Value = "Easy"

$ws.Range("A26").Value = "LeetCode"
$ws.Range("A26").HorizontalAlignment = -4131
$ws.Range("B26").Value = "Stephan"
$ws.Range("B26").HorizontalAlignment = -4131
$ws.Range("C26").Value = "Easy"

# Match the author's final selection/cursor position.
$ws.Range("G25").Select()
